$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (price + 1h volume%), plus a B34/B35 (HuobiToken/Filecoin) row swap
$ws.Range("D2").Value = '27.135.11'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '1.826.67'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("D4").Value = '''1.010'
$ws.Range("E4").Value = '  +0.74%  '
$ws.Range("D5").Value = '''313.19'
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("D7").Value = '''0.4706'
$ws.Range("E7").Value = '  +0.44%  '
$ws.Range("D8").Value = '''0.3653'
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("D9").Value = '''0.07402'
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("D10").Value = '''0.8802'
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("D11").Value = '''20.35'
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").Value = '1.917.95'
$ws.Range("E12").Value = '  +5.10%  '
$ws.Range("D13").Value = '''0.07320'
$ws.Range("E13").Value = '  +2.93%  '
$ws.Range("D14").Value = '''93.33'
$ws.Range("E14").Value = '  +2.17%  '
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '''6.521'
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").Value = '''0.000008709'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = '27.724.94'
$ws.Range("E20").Value = '  +2.81%  '
$ws.Range("D21").Value = '''14.64'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = '''5.240'
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '2.101.81'
$ws.Range("E24").Value = '  +2.67%  '
$ws.Range("D25").Value = '''1.883'
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").Value = '''151.69'
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("D27").Value = '''18.52'
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("D28").Value = '''2.143'
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("D29").Value = '''5.184'
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("D30").Value = '''116.05'
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("D31").Value = '''0.08934'
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").Value = '''1.166'
$ws.Range("E32").Value = '  +0.39%  '
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''2.955'
$ws.Range("E34").Value = '  +1.44%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '''4.510'
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").Value = '''1.009'
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").Value = '''0.05291'
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").Value = '''0.01946'
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("D40").Value = '''2.406'
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("D41").Value = '''2.927'
$ws.Range("E41").Value = '  -1.36%  '
$ws.Range("D42").Value = '''7.221'
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("D43").Value = '''0.5247'
$ws.Range("E43").Value = '  -0.89%  '
$ws.Range("D44").Value = '''0.1642'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '''8.391'
$ws.Range("E45").Value = '  -0.54%  '
$ws.Range("D46").Value = '''0.4867'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").Value = '''10.40'
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("D49").Value = '''104.15'
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("D50").Value = '''1.653'
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").Value = '''0.06298'
$ws.Range("E51").Value = '  -0.04%  '
